# Replace runs of text while preserving the surrounding paragraph/run
# structure: in particular the leading empty <w:r/> run that many
# paragraphs in this document carry before their "real" text run, and
# any direct (explicit) character formatting - bold/italic - applied to
# the run whose text is being edited.
#
# Plain Find/Replace (or a naive Range.Text assignment) makes this
# headless engine coalesce the empty run into the freshly edited run
# whenever the two runs share identical formatting. The fix is to:
#   1. Insert the new text immediately before the old text and then
#      delete the old text (rather than assigning Range.Text directly);
#      this keeps the empty run distinct, but leaves it positioned
#      *after* the new text.
#   2. Cut the freshly-inserted text and paste it back in the very same
#      spot; this forces the engine to re-split it from its neighbour,
#      which restores the original (empty-run-first) ordering.
#   3. Re-apply Bold/Italic explicitly to the new text if (and only if)
#      the original run had *direct* character formatting that differs
#      from its paragraph style's own formatting - otherwise the
#      Insert/Delete/Cut/Paste dance above can lose that direct
#      formatting, and blindly reapplying it would wrongly stamp
#      explicit formatting onto runs that only inherit it from their
#      style (e.g. a Heading1 run, which is bold via its style only).

function Find-ParaIndexByText($d, $oldText) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($oldText)) {
            return $idx
        }
    }
    return -1
}

function Replace-OneParaText($d, $oldText, $newText) {
    $paraIndex = Find-ParaIndexByText $d $oldText
    if ($paraIndex -eq -1) {
        return $false
    }

    $p = $d.Paragraphs($paraIndex)
    $pRange = $p.Range
    $textStart = $pRange.Text.IndexOf($oldText)
    $absStart = $pRange.Start + $textStart
    $absEnd = $absStart + $oldText.Length

    # Figure out whether the text run carries direct Bold/Italic
    # formatting (as opposed to formatting inherited from the
    # paragraph's style), so it can be restored after the edit.
    $srcRng = $d.Range($absStart, $absEnd)
    $styleFont = $p.Style.Font

    $directBold = $false
    $boldVal = $false
    if ($srcRng.Font.Bold -ne $styleFont.Bold) {
        $directBold = $true
        $boldVal = $srcRng.Font.Bold
    }

    $directItalic = $false
    $italicVal = $false
    if ($srcRng.Font.Italic -ne $styleFont.Italic) {
        $directItalic = $true
        $italicVal = $srcRng.Font.Italic
    }

    # Step 1: insert new text before the old text, then delete the old
    # text (keeps any preceding empty run distinct, albeit misordered).
    $insertPoint = $d.Range($absStart, $absStart)
    $insertPoint.InsertBefore($newText)

    $oldStart2 = $absStart + $newText.Length
    $oldEnd2 = $oldStart2 + $oldText.Length
    $delRng = $d.Range($oldStart2, $oldEnd2)
    $delRng.Delete()

    # Restore direct formatting, if any, before the cut/paste step
    # below (Cut/Paste carries a run's own formatting along with it).
    $newRng = $d.Range($absStart, $absStart + $newText.Length)
    if ($directBold) {
        $newRng.Font.Bold = $boldVal
    }
    if ($directItalic) {
        $newRng.Font.Italic = $italicVal
    }

    # Step 2: cut and paste the new text back into the same spot, to
    # restore correct run ordering relative to the empty run.
    $newRng2 = $d.Range($absStart, $absStart + $newText.Length)
    $newRng2.Cut()
    $pasteDest = $d.Range($absStart, $absStart)
    $pasteDest.Paste()

    return $true
}

function Replace-AllParaText($d, $oldText, $newText) {
    $guard = 0
    while ((Replace-OneParaText $d $oldText $newText) -and ($guard -lt 20)) {
        $guard = $guard + 1
    }
}

$d = $word.ActiveDocument

# 1. Title (appears twice: main heading and bold run near the end)
Replace-AllParaText $d `
    "Play Halloween Fortune Free - Review of Exciting Spooky-Themed Game" `
    "Play Halloween Fortune Free: Spooky-themed Slot Game"

# 2. "What we like" bullet list - reordered and reworded
Replace-AllParaText $d `
    "High RTP of 97.6% for frequent substantial payouts" `
    "Easy to play and understand"

Replace-AllParaText $d `
    "Stylish and spooky themed graphics and symbols" `
    "Spooky-themed graphics and detailed symbols"

Replace-AllParaText $d `
    "Win big payouts with the Wild, Scatter, and Bonus symbols" `
    "High RTP of 97.6%"

Replace-AllParaText $d `
    "Available at a range of online casinos" `
    "Plenty of bonus features and opportunities to win big"

# 3. "What we don't like" bullet list - reworded
Replace-AllParaText $d `
    "Not suitable for low-rollers due to the minimum bet amount" `
    "Limited maximum bet of €20"

Replace-AllParaText $d `
    "Free Spins round can be difficult to trigger" `
    "Medium-high volatility may not be suitable for all players"

# 4. Meta description (italic run at the end)
Replace-AllParaText $d `
    "Try your luck with Halloween Fortune, a stylish and thrilling horror-themed slot game. Play free and discover the bonus features for substantial payouts." `
    "Play Halloween Fortune for free and enjoy spooky-themed graphics, high RTP, and exciting bonus features."
